$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 366
$ws.Range("I12").Value = 274
$ws.Range("J12").Value = 550
$ws.Range("K12").Value = 274
$ws.Range("L12").Value = 550
$ws.Range("M12").Value = -104
$ws.Range("N12").Value = -890

$ws.Range("H33").Value = 12890.375
$ws.Range("I33").Value = 16866.334
$ws.Range("J33").Value = 962.5
$ws.Range("K33").Value = 16866.334
$ws.Range("L33").Value = 962.5
$ws.Range("M33").Value = -16637.334
$ws.Range("N33").Value = -1420.5

$ws.Range("H74").Value = 5660.6665
$ws.Range("I74").Value = 5494
$ws.Range("K74").Value = 5494
$ws.Range("M74").Value = -4558

$ws.Range("H76").Value = 4999.5
$ws.Range("I76").Value = 4999.5
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 4999.5
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -4684.5
$ws.Range("N76").ClearContents()

$ws.Range("H77").Value = 5660.6665
$ws.Range("I77").Value = 5494
$ws.Range("K77").Value = 27470
$ws.Range("M77").Value = -22790

$ws.Range("H79").Value = 4999.5
$ws.Range("I79").Value = 4999.5
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 4999.5
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -3907.5
$ws.Range("N79").ClearContents()

$ws.Range("H112").Value = 2310.6365
$ws.Range("I112").Value = 2016.6666
$ws.Range("J112").Value = 2357.0527
$ws.Range("K112").Value = 6049.9998
$ws.Range("L112").Value = 7071.158100000001
$ws.Range("M112").Value = -4941.9998
$ws.Range("N112").Value = -9287.158100000001

$ws.Range("H137").Value = 1277.24
$ws.Range("I137").Value = 921.1429000000001
$ws.Range("K137").Value = 2763.4287
$ws.Range("M137").Value = -213.4287000000004

$ws.Range("H138").Value = 2965.348
$ws.Range("I138").Value = 2163.3076
$ws.Range("J138").Value = 4008
$ws.Range("K138").Value = 6489.9228
$ws.Range("L138").Value = 12024
$ws.Range("M138").Value = -1349.9228
$ws.Range("N138").Value = -22304

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 994.3333
$ws.Range("I4").Value = 994.3333
$ws.Range("K4").Value = 994.3333
$ws.Range("M4").Value = -878.3333

$ws.Range("H32").Value = 18031
$ws.Range("I32").Value = 18707.154
$ws.Range("K32").Value = 18707.154
$ws.Range("M32").Value = -18420.154

$ws.Range("H45").Value = 4903.3
$ws.Range("J45").Value = 6527.5713
$ws.Range("L45").Value = 6527.5713
$ws.Range("N45").Value = -7281.5713

$ws.Range("H61").Value = 4261.4614
$ws.Range("I61").Value = 4261.4614
$ws.Range("K61").Value = 4261.4614
$ws.Range("M61").Value = -4049.4614

$ws.Range("H74").Value = 144677.86
$ws.Range("I74").Value = 168644.5
$ws.Range("J74").Value = 878
$ws.Range("K74").Value = 168644.5
$ws.Range("L74").Value = 878
$ws.Range("M74").Value = -167770.5
$ws.Range("N74").Value = -2626

$ws.Range("H77").Value = 144677.86
$ws.Range("I77").Value = 168644.5
$ws.Range("J77").Value = 878
$ws.Range("K77").Value = 843222.5
$ws.Range("L77").Value = 4390
$ws.Range("M77").Value = -838854.5
$ws.Range("N77").Value = -13126

$ws.Range("H97").Value = 5613.7856
$ws.Range("I97").Value = 7148
$ws.Range("J97").Value = 3242.7273
$ws.Range("K97").Value = 7148
$ws.Range("L97").Value = 3242.7273
$ws.Range("M97").Value = -6652
$ws.Range("N97").Value = -4234.7273

$ws.Range("H132").Value = 26663.512
$ws.Range("I132").Value = 37856.355
$ws.Range("J132").Value = 2555.8462
$ws.Range("K132").Value = 113569.065
$ws.Range("L132").Value = 7667.5386
$ws.Range("M132").Value = -111039.065
$ws.Range("N132").Value = -12727.5386

$ws.Range("H136").Value = 4261.4614
$ws.Range("I136").Value = 4261.4614
$ws.Range("K136").Value = 12784.3842
$ws.Range("M136").Value = -10234.3842

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 252.25
$ws.Range("I22").Value = 565
$ws.Range("J22").Value = 148
$ws.Range("K22").Value = 565
$ws.Range("L22").Value = 148
$ws.Range("M22").Value = -392
$ws.Range("N22").Value = -494

$ws.Range("H81").Value = 214494.83
$ws.Range("J81").Value = 250252
$ws.Range("L81").Value = 250252
$ws.Range("N81").Value = -252374

$ws.Range("H84").Value = 214494.83
$ws.Range("J84").Value = 250252
$ws.Range("L84").Value = 750756
$ws.Range("N84").Value = -761364

$ws.Range("H105").Value = 3284.5
$ws.Range("I105").Value = 3242.4546
$ws.Range("K105").Value = 3242.4546
$ws.Range("M105").Value = -1495.4546

$ws.Range("H134").Value = 2154.32
$ws.Range("I134").Value = 1826.8334
$ws.Range("K134").Value = 5480.5002
$ws.Range("M134").Value = -2945.5002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1079.8
$ws.Range("I16").Value = 1079.8
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1079.8
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -792.8
$ws.Range("N16").ClearContents()

$ws.Range("H17").Value = 9998.5
$ws.Range("I17").Value = 9998.200000000001
$ws.Range("K17").Value = 9998.200000000001
$ws.Range("M17").Value = -9824.200000000001

$ws.Range("H113").Value = 1079.8
$ws.Range("I113").Value = 1079.8
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1079.8
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1090.2
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 319.55554
$ws.Range("I12").Value = 196
$ws.Range("J12").Value = 381.33334
$ws.Range("K12").Value = 588
$ws.Range("L12").Value = 1144.00002
$ws.Range("M12").Value = -415
$ws.Range("N12").Value = -1490.00002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 77.5
$ws.Range("I2").Value = 55
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 55
$ws.Range("L2").Value = 100
$ws.Range("M2").Value = 58
$ws.Range("N2").Value = -326

$ws.Range("H24").Value = 671669
$ws.Range("I24").Value = 1000000
$ws.Range("J24").Value = 15007
$ws.Range("K24").Value = 1000000
$ws.Range("L24").Value = 15007
$ws.Range("M24").Value = -999827
$ws.Range("N24").Value = -15353

$ws.Range("H82").Value = 59998.5
$ws.Range("J82").Value = 59999
$ws.Range("L82").Value = 59999
$ws.Range("N82").Value = -60765

$ws.Range("H85").Value = 59998.5
$ws.Range("J85").Value = 59999
$ws.Range("L85").Value = 59999
$ws.Range("N85").Value = -62651

$ws.Range("H126").Value = 9221.235000000001
$ws.Range("I126").Value = 7532.909
$ws.Range("K126").Value = 22598.727
$ws.Range("M126").Value = -20128.727

$ws.Range("H132").Value = 33732.97
$ws.Range("I132").Value = 39908.54
$ws.Range("K132").Value = 119725.62
$ws.Range("M132").Value = -117195.62

$ws.Range("H135").Value = 50600
$ws.Range("J135").Value = 50600
$ws.Range("L135").Value = 50600
$ws.Range("N135").Value = -60740

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2237.1177
$ws.Range("I22").Value = 648
$ws.Range("J22").Value = 3649.6667
$ws.Range("K22").Value = 648
$ws.Range("L22").Value = 3649.6667
$ws.Range("M22").Value = -353
$ws.Range("N22").Value = -4239.6667

$ws.Range("H27").Value = 2237.1177
$ws.Range("I27").Value = 648
$ws.Range("J27").Value = 3649.6667
$ws.Range("K27").Value = 648
$ws.Range("L27").Value = 3649.6667
$ws.Range("M27").Value = -541
$ws.Range("N27").Value = -3863.6667

$ws.Range("H93").Value = 3047.4443
$ws.Range("I93").Value = 4440
$ws.Range("K93").Value = 4440
$ws.Range("M93").Value = -3192

$ws.Range("H100").Value = 3059.9412
$ws.Range("I100").Value = 2726.9092
$ws.Range("J100").Value = 3670.5
$ws.Range("K100").Value = 2726.9092
$ws.Range("L100").Value = 3670.5
$ws.Range("M100").Value = -2185.9092
$ws.Range("N100").Value = -4752.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 71182.664
$ws.Range("I62").Value = 7782.8
$ws.Range("J62").Value = 150432.5
$ws.Range("K62").Value = 7782.8
$ws.Range("L62").Value = 150432.5
$ws.Range("M62").Value = -7158.8
$ws.Range("N62").Value = -151680.5

$ws.Range("H65").Value = 71182.664
$ws.Range("I65").Value = 7782.8
$ws.Range("J65").Value = 150432.5
$ws.Range("K65").Value = 38914
$ws.Range("L65").Value = 752162.5
$ws.Range("M65").Value = -35794
$ws.Range("N65").Value = -758402.5
